$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "steps counted" (column C) for episodes 1-11 (rows 2-12) ---
$ws.Range("C2").Value = 782
$ws.Range("C3").Value = 486
$ws.Range("C4").Value = 642
$ws.Range("C5").Value = 380
$ws.Range("C6").Value = 576
$ws.Range("C7").Value = 348
$ws.Range("C8").Value = 496
$ws.Range("C9").Value = 736
$ws.Range("C10").Value = 614
$ws.Range("C11").Value = 846
$ws.Range("C12").Value = 748

# Row 12 (episode 11) now reaches the terminal state too: False -> True.
# Use a leading quote so the runtime stores it as literal text "True"
# (matching the workbook's existing inlineStr convention) instead of a
# native boolean, then reset the style so the quote-prefix flag doesn't
# linger in the cell format.
$ws.Range("D12").Value = "'True"
$ws.Range("D12").Style = $ws.Range("C12").Style

# --- Append row 13: episode 12 (ep3 of the 4th experiment) ---
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = 12
$ws.Range("C13").Value = 346
$ws.Range("D13").Value = "'False"
$ws.Range("D13").Style = $ws.Range("C13").Style
